# Update example input file to include Dp / Dp_units columns on the
# "ions" sheet, inserted just before the existing "conc_units" column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ions")

# Insert two new blank columns at I:J; this shifts the old conc_units
# column (and everything in it) from I to K.
$ws.Columns("I:J").Insert()

# Header row
$ws.Range("I1").Value = "Dp"
$ws.Range("J1").Value = "Dp_units"

# Row 2 (CHLORIDE): Dp is a plain number here (no units note in source row)
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = "cm^2/s"

# Rows 3-10: Dp is the same small diffusion-coefficient value, formatted
# in scientific notation like the neighbouring Ds column.
$dpValue = 0.0000019999999999999999
$ws.Range("I3:I10").Value = $dpValue
$ws.Range("I3:I10").NumberFormat = "0.00E+00"
$ws.Range("J3:J10").Value = "cm^2/s"
